$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds two 5-quartile groups (Developed: B:F, Emerging: G:K).
# Processing the factor data into modeling data drops the redundant last
# quartile column from each group, collapsing both groups down to 4
# quartiles (B:E and F:I) and refreshing the computed return values.
# Delete from right to left so earlier column letters stay stable.
$ws.Columns("K").Delete()
$ws.Columns("F").Delete()

# Refresh the forward-return values for the now 8 remaining data columns.
$ws.Range("B4").Value = 0.009906011547455589
$ws.Range("C4").Value = 0.006481457498434323
$ws.Range("D4").Value = 0.006903696708790021
$ws.Range("E4").Value = 0.007273711839882705
$ws.Range("F4").Value = 0.0129296775662297
$ws.Range("G4").Value = 0.008546753634167707
$ws.Range("H4").Value = 0.01097284317135602
$ws.Range("I4").Value = 0.01414604057268972
